# Apply updated "dSF" (column F) values on Sheet1
# Commit message: repull data, push all data, mean calculation

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$changes = @{
    3  = -4
    6  = 1
    7  = 2
    9  = 1
    10 = -2
    11 = -6
    12 = -1
    16 = -1
    17 = 1
    22 = -2
    27 = -1
    28 = -1
    31 = -1
    41 = -1
    42 = -2
    45 = 0
    46 = -1
    49 = 2
    55 = -4
    57 = -6
    60 = 0
    62 = -3
    63 = -6
    67 = -2
    69 = -5
    71 = -3
}

foreach ($row in $changes.Keys) {
    $ws.Range("F$row").Value = $changes[$row]
}

$wb.Save()
